$wb = $excel.ActiveWorkbook

# Bold the header row (row 1) on each data sheet - matches the new bold font
# style (fontId=1) introduced in styles.xml (cellXfs index 2 / 3).
$ws = $wb.Worksheets.Item("phase")
$ws.Range("A1:E1").Font.Bold = $true

$ws = $wb.Worksheets.Item("source")
$ws.Range("A1:E1").Font.Bold = $true

$ws = $wb.Worksheets.Item("prepare")
$ws.Range("A1:E1").Font.Bold = $true

$ws = $wb.Worksheets.Item("train")
$ws.Range("A1:E1").Font.Bold = $true

$ws = $wb.Worksheets.Item("config")
$ws.Range("A1:B1").Font.Bold = $true

# Fix an overlap token during NER training: bump the training iteration
# count in the config sheet from 2 to 20.
$ws.Range("B5").Value = 20

# Make the config sheet the active tab.
$ws.Activate()
